# Add a "lastAccess" column to the Users table in the setup template.
#
# The Users sheet (sheet1) currently has headers: email, role, createdAt
# (A1:C1). We append a new "lastAccess" header in column D, which also
# grows the sheet's used range/dimension from A1:C1 to A1:D1.
#
# Touching the workbook through the COM layer also causes the engine to
# recompute each sheet's <dimension> to match its real used range, which
# incidentally corrects the stale SKUs sheet dimension (A1:G1 -> A1:F1,
# since that sheet's data only spans columns A-F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Cells.Item(1, 4).Value = "lastAccess"
